$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.106.26"
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").Value = "2.301.18"
$ws.Range("E3").Value = "  +0.31%  "

$ws.Range("E4").Value = "  +0.09%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "300.54"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.05%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "97.68"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -1.22%  "

$ws.Range("E7").Value = "  +2.99%  "

$ws.Range("E8").Value = "  +0.07%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.516"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.92%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "36.31"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +0.65%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0792"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.19%  "

$ws.Range("E12").Value = "  +0.66%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "17.73"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -3.68%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "6.88"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.03%  "

$ws.Range("D15").Value = "2.660.74"
$ws.Range("E15").Value = "  +0.38%  "

$ws.Range("D16").Value = "2.311.69"
$ws.Range("E16").Value = "  -1.76%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.789"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -1.31%  "

$ws.Range("D18").Value = "43.010.28"
$ws.Range("E18").Value = "  +0.43%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "13.04"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +3.68%  "

$ws.Range("D20").Value = "0.0₃0911"
$ws.Range("E20").Value = "  +0.88%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.13"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.26%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "68.26"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.87%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "238.13"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +1.15%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.21"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -0.79%  "

$ws.Range("E25").Value = "  -0.40%  "

$ws.Range("E26").Value = "  -0.49%  "

$ws.Range("E27").Value = "  +0.04%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "25.28"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +1.06%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "9.16"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +0.49%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "2.03"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -13.65%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "162.71"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -2.70%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "33.15"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -3.75%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +0.11%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "5.13"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +2.52%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "18.17"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +2.89%  "

$ws.Range("E36").Value = "  +2.13%  "

$ws.Range("E37").Value = "  +0.38%  "

$ws.Range("E38").Value = "  +1.09%  "

$ws.Range("E39").Value = "  +1.15%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.78"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -0.20%  "

$ws.Range("E41").Value = "  +1.18%  "

$ws.Range("E42").Value = "  -1.88%  "

$ws.Range("D43").Value = "2.011.63"
$ws.Range("E43").Value = "  +1.80%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.0287"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -1.33%  "

$ws.Range("E45").Value = "  -6.83%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "10.23"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +1.34%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "17.51"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -0.30%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.85"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -1.27%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "54.36"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -1.41%  "

$ws.Range("D50").Value = "2.535.88"
$ws.Range("E50").Value = "  +0.69%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.53"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -0.43%  "
